$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "Tempo 2"
$ws.Cells.Item(6, 3).Value = 0.0000060000000000000002
$ws.Cells.Item(6, 4).Value = 0.00046900000000000002
$ws.Cells.Item(6, 5).Value = 0.0037200000000000002
$ws.Cells.Item(6, 6).Value = 0.59293799999999997
$ws.Cells.Item(6, 7).Value = 5.6170660000000003

# Row 7 - "Tempo 3"
$ws.Cells.Item(7, 3).Value = 0.0000050000000000000004
$ws.Cells.Item(7, 4).Value = 0.00047199999999999998
$ws.Cells.Item(7, 5).Value = 0.0037100000000000002
$ws.Cells.Item(7, 6).Value = 0.57669400000000004
$ws.Cells.Item(7, 7).Value = 5.7696490000000002

# Row 8 - "Tempo 4"
$ws.Cells.Item(8, 3).Value = 0.0000050000000000000004
$ws.Cells.Item(8, 4).Value = 0.00047399999999999997
$ws.Cells.Item(8, 5).Value = 0.0037299999999999998
$ws.Cells.Item(8, 6).Value = 0.59188099999999999
$ws.Cells.Item(8, 7).Value = 5.8936479999999998

# Row 14 - "Tempo 2"
$ws.Cells.Item(14, 3).Value = 0.0000030000000000000001
$ws.Cells.Item(14, 4).Value = 0.00027500000000000002
$ws.Cells.Item(14, 5).Value = 0.0022430000000000002
$ws.Cells.Item(14, 6).Value = 0.30212499999999998
$ws.Cells.Item(14, 7).Value = 2.5435300000000001
$ws.Cells.Item(14, 8).Value = 1386.16563

# Row 15 - "Tempo 3"
$ws.Cells.Item(15, 3).Value = 0.0000030000000000000001
$ws.Cells.Item(15, 4).Value = 0.00027599999999999999
$ws.Cells.Item(15, 5).Value = 0.002258
$ws.Cells.Item(15, 6).Value = 0.30479699999999998
$ws.Cells.Item(15, 7).Value = 2.5441699999999998
$ws.Cells.Item(15, 8).Value = 1359.10005

# Row 16 - "Tempo 4"
$ws.Cells.Item(16, 3).Value = 0.0000039999999999999998
$ws.Cells.Item(16, 4).Value = 0.00027099999999999997
$ws.Cells.Item(16, 5).Value = 0.0022420000000000001
$ws.Cells.Item(16, 6).Value = 0.30190600000000001
$ws.Cells.Item(16, 7).Value = 2.5185300000000002
$ws.Cells.Item(16, 8).Value = 1389.1965499999999

# Update selection to reflect the active cell when the workbook was last saved
$ws.Range("I20").Select() | Out-Null
